# A new time-log entry was made for row 112 (2014-11-10, 12:45 PM - 2:52 PM,
# 15 min interruption, "Coding" activity), and the cursor moved down to A113
# afterwards - mirroring a user typing a new row into the log and pressing
# Enter/Tab to move on.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A112").Value = 41953              # Date: 11/10/2014
$ws.Range("B112").Value = 0.53125            # Start Time: 12:45 PM
$ws.Range("C112").Value = 0.61944444444444446 # Stop Time: 2:52 PM
$ws.Range("D112").Value = 15                 # Interruption: 15 mins
$ws.Range("F112").Value = "Coding"           # Activity

# Re-assert the Delta formula explicitly (it already lives in the sheet as a
# shared formula anchored at E4) so the new row's cached result is computed
# fresh rather than inherited/stale.
$ws.Range("E112").Formula = '=IF(AND(NOT(ISBLANK(B112)),NOT(ISBLANK(C112))), (C112-B112) * 24 - D112/60, "")'

$excel.CalculateFull()

# Move the selection down to A113, where the user's cursor ended up next.
$ws.Range("A113").Select()

$wb.Save()
